$wb = $excel.ActiveWorkbook

$oldStamp = "February 03 2026 17.29.55 EST"
$newStamp = "February 03 2026 18.05.36 EST"

$aboutSheet = $wb.Worksheets.Item("About")
$dataSheet = $wb.Worksheets.Item("Boundaries and methane sources")

# Update the "About" sheet version lines
$a2 = $aboutSheet.Range("A2")
$a2.Value2 = $a2.Value2.Replace($oldStamp, $newStamp)

$a6 = $aboutSheet.Range("A6")
$a6.Value2 = $a6.Value2.Replace($oldStamp, $newStamp)

# Update the build_version column (S) for rows 2 through 10 on the data sheet
for ($row = 2; $row -le 10; $row++) {
    $cell = $dataSheet.Cells.Item($row, 19)
    $cell.Value2 = $cell.Value2.Replace($oldStamp, $newStamp)
}
